$d = $word.ActiveDocument

# Locate the run that starts the existing "lobster_aov <-" statement; the
# new "lobster_levene" block is inserted immediately before it (still
# inside the same paragraph that currently begins with the
# "#6 run significance tests for lobster size" comment).
$anchor = $d.Content
$found = $anchor.Find.Execute("lobster_aov <-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate insertion anchor 'lobster_aov <-'"
}
$ip = $d.Range($anchor.Start, $anchor.Start)

function Insert-Chunk($range, $text, $styleName) {
    $startPos = $range.Start
    $range.InsertAfter($text)
    $len = $text.Length
    if ($styleName -and $len -gt 0) {
        $styled = $d.Range($startPos, $startPos + $len)
        $styled.Style = $styleName
    }
    $range.SetRange($startPos + $len, $startPos + $len)
}

$VT = [string][char]11
$CR = [string][char]13

# ---- same paragraph: new "lobster_levene <- leveneTest(...)" call ----
Insert-Chunk $ip "lobster_levene <-" "NormalTok"
Insert-Chunk $ip " " "StringTok"
Insert-Chunk $ip "leveneTest" "KeywordTok"
Insert-Chunk $ip "(SIZE" "NormalTok"
Insert-Chunk $ip "~" "OperatorTok"
Insert-Chunk $ip "SITE, " "NormalTok"
Insert-Chunk $ip "data=" "DataTypeTok"
Insert-Chunk $ip "lobster_case_format)" "NormalTok"
Insert-Chunk $ip $VT $null
Insert-Chunk $ip "lobster_levene" "NormalTok"

# ---- new paragraph: the printed Levene's test table ----
Insert-Chunk $ip $CR $null
Insert-Chunk $ip "## Levene's Test for Homogeneity of Variance (center = median)" "VerbatimChar"
Insert-Chunk $ip $VT $null
Insert-Chunk $ip "##         Df F value    Pr(>F)    " "VerbatimChar"
Insert-Chunk $ip $VT $null
Insert-Chunk $ip "## group    4  8.3893 1.065e-06 ***" "VerbatimChar"
Insert-Chunk $ip $VT $null
Insert-Chunk $ip "##       1663                      " "VerbatimChar"
Insert-Chunk $ip $VT $null
Insert-Chunk $ip "## ---" "VerbatimChar"
Insert-Chunk $ip $VT $null
Insert-Chunk $ip "## Signif. codes:  0 '***' 0.001 '**' 0.01 '*' 0.05 '.' 0.1 ' ' 1" "VerbatimChar"

# ---- new paragraph: comment about unequal variances ----
Insert-Chunk $ip $CR $null
Insert-Chunk $ip "#Variences are not equal check to see if another need to do another test beside anova" "CommentTok"
Insert-Chunk $ip $VT $null
Insert-Chunk $ip $VT $null

# The existing "lobster_aov <- aov(...)" run (NormalTok) follows immediately
# after this point, unchanged.
